# "act tablas web jul25" - update Data and Metadata sheets with the July 2025 refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Data sheet: Fecha / Valor time series now goes 2023 down to 1985 (skipping 1989)
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$years = @("2023","2022","2021","2020","2019","2018","2017","2016","2015","2014","2013","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001","2000","1999","1998","1997","1996","1995","1994","1993","1992","1991","1990","1988","1987","1986","1985")
$vals  = @(5.2,5.4,4.7,4.8,5.5,5.8,5.9,6,5.7,5.9,6.1,6,6.3,6.5,6.8,6.3,7,7,6.5,7.1,6.8,6.4,7.1,5.8,6,5.7,6.2,6.3,3.8,4.1,3.8,3.6,2.3,2.9,3.7,4.8,4.7,4.9)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $cellA = $wsData.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $years[$i]
    $cellA.Style = "Normal"

    $wsData.Cells.Item($row, 2).Value = $vals[$i]
}

# ---------------------------------------------------------------------------
# Metadata sheet: blank first key becomes a single space, and a new
# "actualizacion" / "Julio 2025" row is inserted right before "cita".
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Cells.Item(1, 1).Value = " "
$wsMeta.Cells.Item(1, 2).Value = " "

$wsMeta.Rows("9:9").Insert()
$wsMeta.Cells.Item(9, 1).Value = "actualizacion"
$wsMeta.Cells.Item(9, 2).Value = "Julio 2025"
